$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# "lesson 40 topics attached"
# The section "F. Spring" (rows 37-44, lessons 33-40) gets a follow-up block
# of 9 more lesson rows (41-50) appended below it, separated by a blank
# spacer row, mirroring the row/format pattern already used for that
# section. Only the lesson number (col B) and the hour count (col D) are
# filled in for the new rows; topic/date/video columns are left blank
# (to be attached later).
# --------------------------------------------------------------------------

# Copy the cell formatting of the existing "F. Spring" block (B37:F44) down
# to the new block (B46:F53) so the new rows visually match.
$ws.Range("B37:F44").Copy()
$ws.Range("B46").PasteSpecial(-4122)

# Row 44 uses the special "closing row" style (s=22) for column B, but row
# 53 is not the last row of the new block (row 54 is) - reuse the plain
# style (s=24) from one of the interior rows instead.
$ws.Range("B38").Copy()
$ws.Range("B53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Lesson numbers / hour counts for the new rows.
$ws.Range("B46").Value = 41
$ws.Range("D46").Value = 2

$ws.Range("B47").Value = 443
$ws.Range("D47").Value = 2

$ws.Range("B48").Value = 44
$ws.Range("D48").Value = 2

$ws.Range("B49").Value = 45
$ws.Range("D49").Value = 2

$ws.Range("B50").Value = 46
$ws.Range("D50").Value = 2

$ws.Range("B51").Value = 47
$ws.Range("D51").Value = 2

$ws.Range("B52").Value = 48
$ws.Range("D52").Value = 2

$ws.Range("B53").Value = 49
$ws.Range("D53").Value = 2

# Final row of the new block only carries the lesson number, styled like
# the closing row of the previous block (B44).
$ws.Range("B44").Copy()
$ws.Range("B54").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B54").Value = 50

# Row heights: blank spacer row keeps the sheet's old default height,
# the new lesson rows use the (slightly smaller) height now used for
# freshly added rows.
$ws.Rows.Item(45).RowHeight = 14.25
for ($r = 46; $r -le 54; $r++) {
  $ws.Rows.Item($r).RowHeight = 13.8
}

# Leave the selection where the author ended up after typing the new data.
$ws.Range("B55").Select()
